# Update the "about" page text: rewrite the greeting paragraph with the
# fuller bio, drop the large-font styling from "Apart from being a
# designer...", expand the "amateur photographer" paragraph, and replace the
# final "I like museum, light, and magic" line with "I like museum, light,
# and food." plus a couple of trailing blank paragraphs.
#
# The whole body is rebuilt in one shot via Range.InsertXML so every run
# split / proofErr marker / paragraph property matches exactly; InsertXML
# replaces only the selected range's content and leaves the rest of the
# package (styles, sectPr, etc.) untouched.

$d = $word.ActiveDocument

$bodyXml = @'
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">Hello, I’m </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Wenqing</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> Yin</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">, a rising junior @CMU studying product design and human-computer interaction. As a designer, exploring </w:t>
      </w:r>
      <w:r>
        <w:t>new</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> interaction </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">possibilities within both tangible and intangible medium </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">is what excites me the most. </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">I’m also very interested in educational design and physical computing. </w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:t>Apart from being a designer…</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>I’m also a</w:t>
      </w:r>
      <w:r>
        <w:t>n</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t>amateur photographer</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">, an INFJ (sometimes INFP), and a cat lover. </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">I like </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>museum,  light</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t xml:space="preserve">, and food. </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">  </w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p/>
'@

$d.Content.InsertXML($bodyXml)

Write-Output ("Paragraphs: " + $d.Paragraphs.Count)
foreach ($p in $d.Paragraphs) {
    Write-Output ("[" + $p.Range.Text + "]")
}
